$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 85.8724807945396

$ws.Range("N2:N6").Value = $newValue
